$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate() | Out-Null

# Insert a brand-new row 101 (pushes the former row 101 -> row 102, carrying
# its formulas/values with it and re-basing every shared-formula range that
# used to stop at row 101 so it now stops at row 102).
$ws.Rows.Item(101).Insert()

# Copy number/border formatting from the row above (97), which already has
# the C/D/E-only shape we need for the new "Hash / MacBook Air" sample, then
# stamp in the real values/formulas for the new data point (5,000,000).
$ws.Range("A97:E97").Copy()
$ws.Range("A101").PasteSpecial(-4122)

$ws.Range("A101").Value = "MacBook Air"
$ws.Range("B101").Value = "Hash"
$ws.Range("C101").Value = 5000000
$ws.Range("D101").Formula = "=C101^2"
$ws.Range("E101").Formula = "=C101 * LOG(C101, 2)"

# The row-insert left stale placeholder cells behind in the columns (I/K/O)
# that carry shared formulas further down the sheet; row 101 has no data
# there, so drop them.
$ws.Range("F101:O101").ClearContents()

# New column C got wide enough (to fit "MacBook Air"'s row data) that it now
# needs an explicit width.
$ws.Columns.Item(3).ColumnWidth = 8.330729166666666

# Update the view to where the author left off scrolling/selecting.
$ws.Range("P101").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 76
$excel.ActiveWindow.ScrollColumn = 1
